$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "276.74"
Set-TextCell $ws.Range("G2") "16"

Set-TextCell $ws.Range("D3") "21.00"
Set-TextCell $ws.Range("G3") "16"

Set-TextCell $ws.Range("D4") "6.222"
Set-TextCell $ws.Range("G4") "16"

Set-TextCell $ws.Range("D5") "0.06205"
Set-TextCell $ws.Range("G5") "16"

Set-TextCell $ws.Range("D6") "3.578"
Set-TextCell $ws.Range("G6") "16"

Set-TextCell $ws.Range("D7") "6.552"
Set-TextCell $ws.Range("G7") "16"

Set-TextCell $ws.Range("D8") "1.488"
Set-TextCell $ws.Range("G8") "16"

Set-TextCell $ws.Range("D9") "0.8224"
Set-TextCell $ws.Range("G9") "16"

Set-TextCell $ws.Range("D10") "0.1655"
Set-TextCell $ws.Range("G10") "16"

Set-TextCell $ws.Range("D11") "0.08235"
Set-TextCell $ws.Range("G11") "16"

Set-TextCell $ws.Range("D12") "0.03511"
Set-TextCell $ws.Range("G12") "16"

Set-TextCell $ws.Range("D13") "0.03100"
Set-TextCell $ws.Range("G13") "16"

Set-TextCell $ws.Range("D14") "0.09123"
Set-TextCell $ws.Range("G14") "16"

Set-TextCell $ws.Range("D15") "3.773"
Set-TextCell $ws.Range("G15") "16"

Set-TextCell $ws.Range("D16") "0.001608"
Set-TextCell $ws.Range("G16") "16"

Set-TextCell $ws.Range("D17") "0.04691"
Set-TextCell $ws.Range("G17") "16"

Set-TextCell $ws.Range("D18") "0.006464"
Set-TextCell $ws.Range("E18") "17TigerCashTCH"
Set-TextCell $ws.Range("G18") "16"

Set-TextCell $ws.Range("D19") "0.006141"
Set-TextCell $ws.Range("G19") "16"

Set-TextCell $ws.Range("G20") "16"

Set-TextCell $ws.Range("D21") "0.0001500"
Set-TextCell $ws.Range("G21") "16"

Set-TextCell $ws.Range("D22") "3.810"
Set-TextCell $ws.Range("G22") "16"

Set-TextCell $ws.Range("D23") "2.358"
Set-TextCell $ws.Range("G23") "16"

Set-TextCell $ws.Range("D24") "0.01385"
Set-TextCell $ws.Range("G24") "16"

Set-TextCell $ws.Range("D25") "0.3281"
Set-TextCell $ws.Range("G25") "16"

Set-TextCell $ws.Range("G26") "16"

Set-TextCell $ws.Range("G27") "16"

Set-TextCell $ws.Range("G28") "16"

Set-TextCell $ws.Range("G29") "16"

Set-TextCell $ws.Range("G30") "16"

Set-TextCell $ws.Range("G31") "16"

Set-TextCell $ws.Range("G32") "16"

Set-TextCell $ws.Range("G33") "16"

Set-TextCell $ws.Range("G34") "16"

Set-TextCell $ws.Range("G35") "16"

Set-TextCell $ws.Range("G36") "16"

Set-TextCell $ws.Range("G37") "16"

Set-TextCell $ws.Range("G38") "16"

Set-TextCell $ws.Range("G39") "16"

Set-TextCell $ws.Range("G40") "16"

Set-TextCell $ws.Range("D41") "0.007028"
Set-TextCell $ws.Range("G41") "16"

Set-TextCell $ws.Range("B42") "BKEXToken"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws.Range("D42") "0.1104"
Set-TextCell $ws.Range("E42") "41BKEXTokenBKK"
Set-TextCell $ws.Range("G42") "16"

Set-TextCell $ws.Range("B43") "CEJI"
Set-TextCell $ws.Range("C43") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell $ws.Range("D43") "0.003459"
Set-TextCell $ws.Range("E43") "42CEJICEJI"
Set-TextCell $ws.Range("G43") "16"

Set-TextCell $ws.Range("D44") "0.01095"
Set-TextCell $ws.Range("G44") "16"

Set-TextCell $ws.Range("D45") "0.00006273"
Set-TextCell $ws.Range("G45") "16"

Set-TextCell $ws.Range("G46") "16"

Set-TextCell $ws.Range("D47") "0.8458"
Set-TextCell $ws.Range("E47") "46CoinbaseStockTokenCOINBestin24h"
Set-TextCell $ws.Range("G47") "16"

Set-TextCell $ws.Range("D48") "0.001402"
Set-TextCell $ws.Range("G48") "16"

Set-TextCell $ws.Range("G49") "16"

Set-TextCell $ws.Range("G50") "16"

Set-TextCell $ws.Range("G51") "16"
